# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.418.36"
$ws.Range("E2").Value = "  -1.29%  "

# Row 3
$ws.Range("D3").Value = "'1.840.01"
$ws.Range("E3").Value = "  -1.63%  "

# Row 4
$ws.Range("E4").Value = "  -0.08%  "

# Row 5
$ws.Range("D5").Value = "'261.68"
$ws.Range("E5").Value = "  -5.58%  "

# Row 6
$ws.Range("E6").Value = "  -0.05%  "

# Row 7
$ws.Range("D7").Value = "'0.5197"
$ws.Range("E7").Value = "  -1.56%  "

# Row 8
$ws.Range("D8").Value = "'0.3273"
$ws.Range("E8").Value = "  -4.24%  "

# Row 9
$ws.Range("D9").Value = "'0.06783"
$ws.Range("E9").Value = "  -2.16%  "

# Row 10
$ws.Range("E10").Value = "  -6.60%  "

# Row 11
$ws.Range("D11").Value = "'0.7737"
$ws.Range("E11").Value = "  -3.60%  "

# Row 12
$ws.Range("D12").Value = "'0.07749"
$ws.Range("E12").Value = "  -0.23%  "

# Row 13
$ws.Range("D13").Value = "'1.834.36"
$ws.Range("E13").Value = "  -2.06%  "

# Row 14
$ws.Range("D14").Value = "'87.92"
$ws.Range("E14").Value = "  -2.56%  "

# Row 15
$ws.Range("D15").Value = "'5.004"
$ws.Range("E15").Value = "  -3.31%  "

# Row 16
$ws.Range("E16").Value = "  -0.19%  "

# Row 17
$ws.Range("D17").Value = "'13.92"
$ws.Range("E17").Value = "  -4.37%  "

# Row 18
$ws.Range("E18").Value = "  -0.07%  "

# Row 19
$ws.Range("D19").Value = "'0.000007933"
$ws.Range("E19").Value = "  -1.23%  "

# Row 20
$ws.Range("D20").Value = "'26.483.58"
$ws.Range("E20").Value = "  -1.18%  "

# Row 21
$ws.Range("D21").Value = "'2.076.64"
$ws.Range("E21").Value = "  -0.83%  "

# Row 22
$ws.Range("D22").Value = "'4.608"
$ws.Range("E22").Value = "  -2.89%  "

# Row 23
$ws.Range("D23").Value = "'9.548"
$ws.Range("E23").Value = "  -4.82%  "

# Row 24
$ws.Range("D24").Value = "'6.011"
$ws.Range("E24").Value = "  -2.42%  "

# Row 25
$ws.Range("D25").Value = "'145.53"
$ws.Range("E25").Value = "  -0.56%  "

# Row 26
$ws.Range("D26").Value = "'2.192"
$ws.Range("E26").Value = "  -7.25%  "

# Row 27
$ws.Range("E27").Value = "  -0.14%  "

# Row 28
$ws.Range("E28").Value = "  -2.22%  "

# Row 29
$ws.Range("D29").Value = "'111.58"
$ws.Range("E29").Value = "  -1.60%  "

# Row 30
$ws.Range("D30").Value = "'4.198"
$ws.Range("E30").Value = "  -3.16%  "

# Row 31
$ws.Range("D31").Value = "'4.123"
$ws.Range("E31").Value = "  -4.64%  "

# Row 32
$ws.Range("D32").Value = "'0.08698"
$ws.Range("E32").Value = "  -2.41%  "

# Row 33
$ws.Range("D33").Value = "'0.04820"
$ws.Range("E33").Value = "  -2.26%  "

# Row 34
$ws.Range("E34").Value = "  -2.90%  "

# Row 35
$ws.Range("D35").Value = "'0.7178"
$ws.Range("E35").Value = "  -1.48%  "

# Row 36
$ws.Range("E36").Value = "  -1.29%  "

# Row 37
$ws.Range("D37").Value = "'3.088"
$ws.Range("E37").Value = "  -5.55%  "

# Row 38
$ws.Range("D38").Value = "'0.01781"
$ws.Range("E38").Value = "  -3.91%  "

# Row 39
$ws.Range("D39").Value = "'2.228"
$ws.Range("E39").Value = "  -3.97%  "

# Row 40
$ws.Range("D40").Value = "'0.4839"
$ws.Range("E40").Value = "  -5.76%  "

# Row 41
$ws.Range("D41").Value = "'112.40"
$ws.Range("E41").Value = "  -3.11%  "

# Row 42
$ws.Range("D42").Value = "'0.8998"
$ws.Range("E42").Value = "  -5.18%  "

# Row 43
$ws.Range("D43").Value = "'6.079"
$ws.Range("E43").Value = "  -1.28%  "

# Row 44
$ws.Range("D44").Value = "'1.000"
$ws.Range("E44").Value = "  -0.05%  "

# Row 45
$ws.Range("D45").Value = "'7.728"
$ws.Range("E45").Value = "  -4.39%  "

# Row 46
$ws.Range("B46").Value = "Decentraland"
$ws.Range("C46").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D46").Value = "'0.4152"
$ws.Range("E46").Value = "  -6.94%  "

# Row 47
$ws.Range("B47").Value = "Cronos"
$ws.Range("C47").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D47").Value = "'0.05901"
$ws.Range("E47").Value = "  -2.45%  "

# Row 48
$ws.Range("D48").Value = "'9.004"
$ws.Range("E48").Value = "  -2.75%  "

# Row 49
$ws.Range("D49").Value = "'35.08"
$ws.Range("E49").Value = "  -3.40%  "

# Row 50
$ws.Range("D50").Value = "'0.1219"
$ws.Range("E50").Value = "  -8.96%  "

# Row 51
$ws.Range("D51").Value = "'0.8855"
$ws.Range("E51").Value = "  +0.04%  "
